# Add two new columns (I, J) with headers "I0" and "IF", matching the
# existing header formatting used by column H, plus their row-2 data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy H1's formatting (bold font, thin border box, centered/top-aligned)
# onto the new header cells so they reuse the same cell style as the rest
# of the header row instead of getting ad-hoc formatting.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data cells for row 2
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
